# "excel writer date impl" - the writer now emits the formatted date/time
# string as literal text (instead of a numeric date serial) for column F,
# using a custom "d/m/yyyy h:mm" number format on those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (date column) - replace the numeric date serials with their
# textual "d/m/yyyy h:mm" representation. Values are assigned in the same
# order the source writer produced its string table (alphabetical on the
# formatted text) so the resulting shared-string table lines up with it.
$ws.Range("F2").Value  = "13/08/2017 14:37"
$ws.Range("F7").Value  = "13/08/2017 14:42"
$ws.Range("F8").Value  = "13/08/2017 14:43"
$ws.Range("F9").Value  = "13/08/2017 14:44"
$ws.Range("F11").Value = "13/08/2017 14:46"
$ws.Range("F3").Value  = "14/08/2017 14:38"
$ws.Range("F4").Value  = "15/08/2017 14:39"
$ws.Range("F5").Value  = "16/08/2017 14:40"
$ws.Range("F6").Value  = "17/08/2017 14:41"
$ws.Range("F10").Value = "22/08/2017 14:45"

# Give the date column its own custom number format (replaces the built-in
# date format previously applied to those cells).
$ws.Range("F2:F11").NumberFormat = "d/m/yyyy\ h:mm"

# Move the active selection.
$null = $ws.Range("I10").Select()

# Print setup for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
